{"js": "// Update the three-digit-by-one-digit multiplication answers in the table.\n// Each old equation string is unique in the document, so we can safely\n// search for the exact old text and replace it with the new text.\nconst replacements = [\n  [\"770\u00d72=1540\", \"424\u00d72=848\"],\n  [\"133\u00d76=798\", \"811\u00d76=4866\"],\n  [\"837\u00d79=7533\", \"748\u00d79=6732\"],\n  [\"790\u00d75=3950\", \"128\u00d72=256\"],\n  [\"282\u00d77=1974\", \"127\u00d76=762\"],\n  [\"896\u00d75=4480\", \"235\u00d75=1175\"],\n  [\"357\u00d79=3213\", \"563\u00d74=2252\"],\n  [\"886\u00d73=2658\", \"491\u00d75=2455\"],\n  [\"962\u00d75=4810\", \"362\u00d76=2172\"],\n  [\"878\u00d78=7024\", \"984\u00d79=8856\"],\n  [\"967\u00d78=7736\", \"563\u00d76=3378\"],\n  [\"466\u00d76=2796\", \"313\u00d77=2191\"],\n  [\"487\u00d75=2435\", \"228\u00d79=2052\"],\n  [\"162\u00d79=1458\", \"769\u00d78=6152\"],\n  [\"553\u00d78=4424\", \"289\u00d79=2601\"],\n  [\"286\u00d72=572\", \"904\u00d73=2712\"],\n  [\"219\u00d72=438\", \"362\u00d76=2172\"],\n  [\"731\u00d77=5117\", \"233\u00d79=2097\"],\n  [\"747\u00d72=1494\", \"246\u00d73=738\"],\n  [\"216\u00d74=864\", \"837\u00d73=2511\"],\n  [\"174\u00d78=1392\", \"806\u00d79=7254\"],\n  [\"136\u00d73=408\", \"178\u00d74=712\"],\n  [\"178\u00d72=356\", \"555\u00d73=1665\"],\n  [\"888\u00d75=4440\", \"919\u00d72=1838\"],\n  [\"884\u00d79=7956\", \"559\u00d76=3354\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit-by-one-digit multiplication answers in the table.\n# Each old equation string is unique in the document, so a simple\n# Find/Replace (match case, whole document) per pair is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"770\u00d72=1540\", \"424\u00d72=848\"),\n    @(\"133\u00d76=798\", \"811\u00d76=4866\"),\n    @(\"837\u00d79=7533\", \"748\u00d79=6732\"),\n    @(\"790\u00d75=3950\", \"128\u00d72=256\"),\n    @(\"282\u00d77=1974\", \"127\u00d76=762\"),\n    @(\"896\u00d75=4480\", \"235\u00d75=1175\"),\n    @(\"357\u00d79=3213\", \"563\u00d74=2252\"),\n    @(\"886\u00d73=2658\", \"491\u00d75=2455\"),\n    @(\"962\u00d75=4810\", \"362\u00d76=2172\"),\n    @(\"878\u00d78=7024\", \"984\u00d79=8856\"),\n    @(\"967\u00d78=7736\", \"563\u00d76=3378\"),\n    @(\"466\u00d76=2796\", \"313\u00d77=2191\"),\n    @(\"487\u00d75=2435\", \"228\u00d79=2052\"),\n    @(\"162\u00d79=1458\", \"769\u00d78=6152\"),\n    @(\"553\u00d78=4424\", \"289\u00d79=2601\"),\n    @(\"286\u00d72=572\", \"904\u00d73=2712\"),\n    @(\"219\u00d72=438\", \"362\u00d76=2172\"),\n    @(\"731\u00d77=5117\", \"233\u00d79=2097\"),\n    @(\"747\u00d72=1494\", \"246\u00d73=738\"),\n    @(\"216\u00d74=864\", \"837\u00d73=2511\"),\n    @(\"174\u00d78=1392\", \"806\u00d79=7254\"),\n    @(\"136\u00d73=408\", \"178\u00d74=712\"),\n    @(\"178\u00d72=356\", \"555\u00d73=1665\"),\n    @(\"888\u00d75=4440\", \"919\u00d72=1838\"),\n    @(\"884\u00d79=7956\", \"559\u00d76=3354\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
